$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 to make room for the new draw (15 July 2017),
# shifting all existing draw rows down by one.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row with the new draw's data.
$ws.Range("A6").Value = "15 July 2017"
$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 17
$ws.Range("D6").Value = 26
$ws.Range("E6").Value = 30
$ws.Range("F6").Value = 46
$ws.Range("G6").Value = 48
$ws.Range("H6").Value = 34
